$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three remaining "attendee" assets to show they are finished/original
$ws.Range("D8").Value = "final"
$ws.Range("E8").Value = "original"

$ws.Range("D10").Value = "final"
$ws.Range("E10").Value = "original"

$ws.Range("D11").Value = "final"
$ws.Range("E11").Value = "original"

# Update selected cell to reflect the last reviewed row
$ws.Range("E12").Select()
